# The deck ships two theme parts:
#   theme/theme1.xml -> clrScheme "Office"   (name="Office Theme")
#   theme/theme2.xml -> clrScheme "Integral" (name="Integral"), linked from slideMaster1
# The target edit swaps the two themes' contents: theme1 becomes "Integral"
# and theme2 becomes "Office Theme". The live/used theme (the one driving
# every slide through slideMaster1) is theme2.xml, reachable here through
# the ThemeColorScheme on any slide. Re-point its 12 scheme colors to the
# classic Office palette (the colors theme1 currently holds).

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index order (VBA ppThemeColorIndex): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $b = $hex -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $r = ($hex -shr 16) -band 0xFF
    $vbaRgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $vbaRgb
}
